# Updates the cryptos list (prices / 1h volume %) per the latest scrape,
# including a swap of the Maker / ApeXProtocol rows (47 and 48).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "50.931.01"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.909.09"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'368.54"
$ws.Range("E5").Value = "  +5.24%  "
$ws.Range("D6").Value = "'102.57"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("D7").Value = "'0.540"
$ws.Range("E7").Value = "  -2.55%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.582"
$ws.Range("D10").Value = "'36.79"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'0.0832"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").Value = "3.358.16"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "'7.36"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "2.899.00"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'0.924"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").Value = "50.859.27"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'3.21"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").Value = "'7.17"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("E21").Value = "  -3.97%  "
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").Value = "'67.96"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "'258.02"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "'2.66"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("D29").Value = "'25.53"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").Value = "'7.08"
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").Value = "'6.27"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").Value = "'9.86"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").Value = "'51.20"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").Value = "'34.21"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").Value = "'17.01"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "'2.60"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'1.84"
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "'119.22"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.016.91"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'2.31"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("D50").Value = "3.185.31"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  +0.30%  "
